$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the first 5 data rows (rows 2-6) with new values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 160

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 155

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 149

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 126

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 125

# Delete rows 7-11 which are no longer needed
$ws.Range("A7:B11").EntireRow.Delete()
